$wb = $excel.ActiveWorkbook

# --- Dash_1 sheet: remove the daily rows below the first data day (rows 5-31), ---
# --- keep the freeze pane / header, and move the active selection back to row 1 ---
$wsDash = $wb.Worksheets.Item("Dash_1")
$wsDash.Rows("5:31").Delete()

# --- Exhibition sheet: the last header cell (I1) used a redundant "apply fill" ---
# --- border style identical to the rest of the header row (H1); normalize it ---
$wsExhibition = $wb.Worksheets.Item("Exhibition")
$wsExhibition.Range("H1").Copy()
$wsExhibition.Range("I1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Active tab / selection moves from Exhibition back to Dash_1 ---
$wsDash.Activate()
$wsDash.Range("L1").Select()
